$wb = $excel.ActiveWorkbook

# Sheet "展览" - column F ("想去人数") updates: refresh generated attendance counts
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Cells.Item(2, 6).Value = 880
$wsExhibit.Cells.Item(5, 6).Value = 1183
$wsExhibit.Cells.Item(7, 6).Value = 4332
$wsExhibit.Cells.Item(8, 6).Value = 2589
$wsExhibit.Cells.Item(10, 6).Value = 2494
$wsExhibit.Cells.Item(14, 6).Value = 1653
$wsExhibit.Cells.Item(15, 6).Value = 657
$wsExhibit.Cells.Item(16, 6).Value = 18
$wsExhibit.Cells.Item(17, 6).Value = 109
$wsExhibit.Cells.Item(18, 6).Value = 317
$wsExhibit.Cells.Item(19, 6).Value = 25
$wsExhibit.Cells.Item(20, 6).Value = 270
$wsExhibit.Cells.Item(21, 6).Value = 73
$wsExhibit.Cells.Item(22, 6).Value = 22
$wsExhibit.Cells.Item(23, 6).Value = 472
$wsExhibit.Cells.Item(25, 6).Value = 87
$wsExhibit.Cells.Item(26, 6).Value = 532
$wsExhibit.Cells.Item(27, 6).Value = 689
$wsExhibit.Cells.Item(28, 6).Value = 101
$wsExhibit.Cells.Item(29, 6).Value = 79
$wsExhibit.Cells.Item(30, 6).Value = 397
$wsExhibit.Cells.Item(31, 6).Value = 46
$wsExhibit.Cells.Item(32, 6).Value = 1614
$wsExhibit.Cells.Item(33, 6).Value = 996
$wsExhibit.Cells.Item(34, 6).Value = 108
$wsExhibit.Cells.Item(35, 6).Value = 17
$wsExhibit.Cells.Item(36, 6).Value = 1094
$wsExhibit.Cells.Item(37, 6).Value = 2027
$wsExhibit.Cells.Item(38, 6).Value = 258
$wsExhibit.Cells.Item(40, 6).Value = 540
$wsExhibit.Cells.Item(42, 6).Value = 22
$wsExhibit.Cells.Item(43, 6).Value = 647
$wsExhibit.Cells.Item(44, 6).Value = 1305
$wsExhibit.Cells.Item(45, 6).Value = 86
$wsExhibit.Cells.Item(47, 6).Value = 427

# Sheet "全部类型" - column F ("想去人数") updates: refresh generated attendance counts
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Cells.Item(2, 6).Value = 880
$wsAll.Cells.Item(3, 6).Value = 1183
$wsAll.Cells.Item(6, 6).Value = 4332
$wsAll.Cells.Item(7, 6).Value = 2589
$wsAll.Cells.Item(8, 6).Value = 2494
$wsAll.Cells.Item(9, 6).Value = 1653
$wsAll.Cells.Item(12, 6).Value = 657
$wsAll.Cells.Item(13, 6).Value = 18
$wsAll.Cells.Item(14, 6).Value = 109
$wsAll.Cells.Item(15, 6).Value = 317
$wsAll.Cells.Item(16, 6).Value = 25
$wsAll.Cells.Item(17, 6).Value = 270
$wsAll.Cells.Item(18, 6).Value = 73
$wsAll.Cells.Item(19, 6).Value = 472
$wsAll.Cells.Item(21, 6).Value = 87
$wsAll.Cells.Item(22, 6).Value = 532
$wsAll.Cells.Item(23, 6).Value = 689
$wsAll.Cells.Item(24, 6).Value = 101
$wsAll.Cells.Item(28, 6).Value = 79
$wsAll.Cells.Item(29, 6).Value = 397
$wsAll.Cells.Item(30, 6).Value = 1614
$wsAll.Cells.Item(31, 6).Value = 996
$wsAll.Cells.Item(32, 6).Value = 108
$wsAll.Cells.Item(34, 6).Value = 2028
$wsAll.Cells.Item(35, 6).Value = 258
$wsAll.Cells.Item(40, 6).Value = 540
$wsAll.Cells.Item(42, 6).Value = 22
$wsAll.Cells.Item(43, 6).Value = 647
$wsAll.Cells.Item(44, 6).Value = 1305
$wsAll.Cells.Item(46, 6).Value = 86
$wsAll.Cells.Item(47, 6).Value = 427

Write-Host "edits applied"
